$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the existing "actions" table (server-->avr) with one more entry
$ws.Range("B6").Value = "0x03"
$ws.Range("C6").Value = "battery level callback"

# New "avr-->server" table
$ws.Range("A8").Value = "avr-->server"
$ws.Range("B8").Value = "0x01"
$ws.Range("C8").Value = "ready for transmission"

$ws.Range("B9").Value = "0x02"
$ws.Range("C9").Value = "battery level transmission (2 bytes)"

$ws.Range("B10").Value = "0x03"
$ws.Range("C10").Value = "bad action"

$ws.Range("B11").Value = "0x04"
$ws.Range("C11").Value = "bad key"

$ws.Range("B12").Value = "0x05"
$ws.Range("C12").Value = "OK"

# Widen column C to fit the new, longer text
$ws.Columns("C").ColumnWidth = 33.65

# Move the active selection, matching the author's final cursor position
$ws.Range("D7").Select()
